$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark that currently sits
#    between "multiple" and " GUI testing framework tools cu".
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the paragraph that ends with
#    "their strengths and weaknesses as a whole." (end of the
#    "Mission" section) so we can re-add the bookmark there once the
#    trailing "Executive Summary" paragraphs are removed.
# ------------------------------------------------------------------
$missionBodyParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pText = $d.Paragraphs.Item($i).Range.Text
    if ($pText.Contains("their strengths and weaknesses as a whole.")) {
        $missionBodyParaIndex = $i
        break
    }
}

$missionPara = $d.Paragraphs.Item($missionBodyParaIndex)
$endPos = $missionPara.Range.End - 1   # position right before the paragraph mark

# ------------------------------------------------------------------
# 3. Delete the "Executive Summary" heading paragraph and the
#    following "TBD - add summary..." paragraph - they are the two
#    paragraphs immediately after the Mission body paragraph and run
#    to the end of the document body content.
# ------------------------------------------------------------------
$execSummaryPara = $d.Paragraphs.Item($missionBodyParaIndex + 1)
$deleteRange = $d.Range($execSummaryPara.Range.Start, $d.Content.End)
$deleteRange.Delete()

# ------------------------------------------------------------------
# 4. Re-add the "_GoBack" bookmark at the end of the Mission body
#    paragraph text (right before its paragraph mark).
#
#    The runtime's Bookmarks.Add mishandles a zero-length range that
#    sits exactly at "end of paragraph text" (i.e. immediately before
#    a paragraph mark) - it silently snaps the bookmark elsewhere.
#    Work around this by temporarily inserting a placeholder
#    character after the target position (which moves the position
#    away from the paragraph boundary), adding the bookmark there,
#    and then removing the placeholder again. The collapsed bookmark
#    remains correctly anchored once the placeholder is deleted.
# ------------------------------------------------------------------
$placeholderRange = $d.Range($endPos, $endPos)
$placeholderRange.InsertAfter("X")

$bookmarkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholderCharRange = $d.Range($endPos, $endPos + 1)
$placeholderCharRange.Delete()
